$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D gets manually widened (no longer auto "best fit") ---
$ws.Columns("D").ColumnWidth = 31

# --- Work in manual calculation mode so we can capture the same kind of
#     "stale formula cache" the authored workbook ended up with. ---
$excel.Calculation = -4135   # xlCalculationManual

# Row 2 / Row 3: the "Output 1" keyword cell (I2) is repurposed to hold the
# literal query result "56" instead of the descriptive keyword text. D3 is
# a live formula "=I2" whose cached display text is captured as "out" (the
# in-between value I2 held when the sheet was last recalculated).
$ws.Range("I2").Value = "out"
$excel.Calculate()
$ws.Range("I2").Value = "'56"
$ws.Range("I2").ClearFormats()

# Row 3: the expected check value becomes 56 (matches the new query result).
$ws.Range("E3").Value = 56

# Row 4 / Row 5: the second "Output 1" keyword cell (I4) is repurposed to
# hold the literal ky_ba value "9041383009". D5's cached formula text is
# left stale (it keeps showing the previous "output ky_ba" text).
$ws.Range("I4").Value = "'9041383009"
$ws.Range("I4").ClearFormats()

# --- Selection / scroll position at save time ---
$ws.Activate()
$ws.Range("I3").Select()
